$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) — sheet index 1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3880
$ws1.Range("F4").Value = 2305
$ws1.Range("F11").Value = 1436
$ws1.Range("F13").Value = 2553
$ws1.Range("F14").Value = 180

# Sheet "全部类型" (All types) — sheet index 4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3880
$ws4.Range("F4").Value = 2305
$ws4.Range("F14").Value = 1436
$ws4.Range("F16").Value = 2553
$ws4.Range("F17").Value = 180
